$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Yht" (totals) row 22 is pushed down to row 23, and a new
# work-log entry (previously only implied by the totals formula range)
# is inserted as the new row 22.

# 1) Copy the formatting of an existing data row (row 11 uses the same
#    date / hours / wrapped-description layout we need) onto row 22.
$ws.Range("B11:D11").Copy() | Out-Null
$ws.Range("B22:D22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B22").Value = 45340
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = "Tein lievästi tyylien kanssa säätöä. Hankin 3 testaajaa ja sain heiltä ongelmia. Korjasin kyseiset bugit, aloin suunitelemaan ikonia."
$ws.Rows.Item(22).RowHeight = 56.25

# 2) Recreate the totals row one line further down, reusing the
#    formatting of the header row (same "Good" style as the old row 22).
$ws.Range("B5:D5").Copy() | Out-Null
$ws.Range("B23:D23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B23").Value = "Yht"
$ws.Range("C23").Formula = "=SUM(C6:C22)"
$ws.Range("D23").Value = $null
$ws.Rows.Item(23).RowHeight = 18.75

$excel.CutCopyMode = 0

# 3) Update the selection / active cell to mirror the new layout.
$ws.Range("D24").Select()

$wb.Save()
